# Generate Report for Handback
# Updates the "26285073-2a95-49bb-b4e6-ac4255bd9f16.md" row (row 7) on both the
# "zh-cn" and "de-de" status sheets: the handback for that file has now come
# back, but its source version is stale, so we record the latest target file,
# the handback datetime, a hyperlink to the target file, and the error detail
# explaining the version mismatch.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b91ed28b5f5746569a3499da9c911ce4dbabfce1/e2e/26285073-2a95-49bb-b4e6-ac4255bd9f16.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e6de207cd881ec3d0c0f765fbe1ced377f67b761/e2e/26285073-2a95-49bb-b4e6-ac4255bd9f16.md."

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$targetFileNameZhCn = $wsZhCn.Range("A7").Text
$latestHandoffFileZhCn = $wsZhCn.Range("G7").Text

# Latest Target File (I7) - same display name as the source file, turned
# into a hyperlink like the other rows in the sheet.
$wsZhCn.Range("I7").Value = $targetFileNameZhCn
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/837253e5b7cdc730ea6dc5803c3379d487830f61/e2e/26285073-2a95-49bb-b4e6-ac4255bd9f16.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetFileNameZhCn) | Out-Null

# Latest Handback File (J7) - same xlf as the latest handoff file (G7).
$wsZhCn.Range("J7").Value = $latestHandoffFileZhCn

# Latest Handback DateTime (K7).
$wsZhCn.Range("K7").Value = "2016-09-03 20:57:46"

# Error Detail (P7) - version mismatch message.
$wsZhCn.Range("P7").Value = $errorDetail

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$targetFileNameDeDe = $wsDeDe.Range("A7").Text
$latestHandoffFileDeDe = $wsDeDe.Range("G7").Text

# Latest Target File (I7) - same display name as the source file, turned
# into a hyperlink like the other rows in the sheet.
$wsDeDe.Range("I7").Value = $targetFileNameDeDe
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/134ba3c1ebe677b72668722d1d6846e6305220bc/e2e/26285073-2a95-49bb-b4e6-ac4255bd9f16.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $targetFileNameDeDe) | Out-Null

# Latest Handback File (J7) - same xlf as the latest handoff file (G7).
$wsDeDe.Range("J7").Value = $latestHandoffFileDeDe

# Latest Handback DateTime (K7).
$wsDeDe.Range("K7").Value = "2016-09-03 20:57:53"

# Error Detail (P7) - version mismatch message.
$wsDeDe.Range("P7").Value = $errorDetail
